$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cells value as literal text (force Text number format so
# numeric-looking strings like "51.468.45" or "1.00" are not auto-converted
# to numbers / floats, matching the original inline-string cell content).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Rows 47 and 48 swap coin data (TheGraph <-> ApeXProtocol) with updated price/volume values
Set-TextValue $ws.Range("B47") 'ApeXProtocol'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D47") '2.36'
Set-TextValue $ws.Range("E47") '  +1.14%  '

Set-TextValue $ws.Range("B48") 'TheGraph'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range("D48") '0.270'
Set-TextValue $ws.Range("E48") '  -0.52%  '

# Price / volume updates for the other rows
Set-TextValue $ws.Range("D2") '51.468.45'
Set-TextValue $ws.Range("E2") '  +0.86%  '
Set-TextValue $ws.Range("D3") '2.980.79'
Set-TextValue $ws.Range("E3") '  +1.12%  '
Set-TextValue $ws.Range("E4") '  -0.01%  '
Set-TextValue $ws.Range("D5") '381.28'
Set-TextValue $ws.Range("E5") '  +1.38%  '
Set-TextValue $ws.Range("D6") '103.84'
Set-TextValue $ws.Range("E6") '  +2.48%  '
Set-TextValue $ws.Range("E7") '  +0.68%  '
Set-TextValue $ws.Range("E8") '  -0.01%  '
Set-TextValue $ws.Range("D9") '0.590'
Set-TextValue $ws.Range("E9") '  -0.09%  '
Set-TextValue $ws.Range("D10") '36.59'
Set-TextValue $ws.Range("E10") '  +0.56%  '
Set-TextValue $ws.Range("E11") '  -0.96%  '
Set-TextValue $ws.Range("D12") '0.0857'
Set-TextValue $ws.Range("E12") '  +0.66%  '
Set-TextValue $ws.Range("D13") '3.454.66'
Set-TextValue $ws.Range("E13") '  +1.51%  '
Set-TextValue $ws.Range("D14") '7.84'
Set-TextValue $ws.Range("E14") '  +2.70%  '
Set-TextValue $ws.Range("D15") '18.45'
Set-TextValue $ws.Range("D16") '2.985.29'
Set-TextValue $ws.Range("E16") '  +1.40%  '
Set-TextValue $ws.Range("D17") '11.21'
Set-TextValue $ws.Range("E17") '  -0.68%  '
Set-TextValue $ws.Range("D18") '0.995'
Set-TextValue $ws.Range("E18") '  -0.27%  '
Set-TextValue $ws.Range("D19") '51.490.08'
Set-TextValue $ws.Range("E19") '  +0.96%  '
Set-TextValue $ws.Range("E20") '  +0.26%  '
Set-TextValue $ws.Range("D21") '12.56'
Set-TextValue $ws.Range("E21") '  +0.45%  '
Set-TextValue $ws.Range("D22") '0.0₃0961'
Set-TextValue $ws.Range("E22") '  +0.53%  '
Set-TextValue $ws.Range("D23") '70.25'
Set-TextValue $ws.Range("E23") '  +1.93%  '
Set-TextValue $ws.Range("D24") '266.83'
Set-TextValue $ws.Range("E24") '  +0.11%  '
Set-TextValue $ws.Range("D25") '3.22'
Set-TextValue $ws.Range("E25") '  +1.51%  '
Set-TextValue $ws.Range("D26") '7.82'
Set-TextValue $ws.Range("E26") '  -4.70%  '
Set-TextValue $ws.Range("E27") '  -2.95%  '
Set-TextValue $ws.Range("D28") '0.169'
Set-TextValue $ws.Range("E28") '  +3.69%  '
Set-TextValue $ws.Range("D29") '1.00'
Set-TextValue $ws.Range("E29") '  +0.00%  '
Set-TextValue $ws.Range("E30") '  +1.30%  '
Set-TextValue $ws.Range("E31") '  -1.21%  '
Set-TextValue $ws.Range("E32") '  +3.92%  '
Set-TextValue $ws.Range("D33") '34.56'
Set-TextValue $ws.Range("E33") '  +3.35%  '
Set-TextValue $ws.Range("D34") '51.34'
Set-TextValue $ws.Range("E34") '  +0.74%  '
Set-TextValue $ws.Range("E35") '  +0.25%  '
Set-TextValue $ws.Range("E36") '  +0.91%  '
Set-TextValue $ws.Range("E37") '  +0.09%  '
Set-TextValue $ws.Range("D38") '3.27'
Set-TextValue $ws.Range("E38") '  +3.86%  '
Set-TextValue $ws.Range("E39") '  +2.50%  '
Set-TextValue $ws.Range("E40") '  +3.88%  '
Set-TextValue $ws.Range("E41") '  +0.65%  '
Set-TextValue $ws.Range("E42") '  +1.77%  '
Set-TextValue $ws.Range("E43") '  +12.30%  '
Set-TextValue $ws.Range("D44") '126.02'
Set-TextValue $ws.Range("E44") '  +5.37%  '
Set-TextValue $ws.Range("D45") '21.32'
Set-TextValue $ws.Range("E45") '  +0.27%  '
Set-TextValue $ws.Range("E46") '  -0.05%  '
Set-TextValue $ws.Range("D49") '2.023.34'
Set-TextValue $ws.Range("E49") '  +1.44%  '
Set-TextValue $ws.Range("D50") '3.278.30'
Set-TextValue $ws.Range("E50") '  +1.22%  '
Set-TextValue $ws.Range("E51") '  +0.05%  '
